$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 9 (shifts everything from row 9 down by one)
$ws.Rows.Item(9).Insert()

# New label cell (A9) - reuse the same style as the other label cells in this
# section (A4:A8), i.e. style index 11 -> mimic via existing cell's style
$ws.Range("A9").Value = "OrcID"
$ws.Range("A9").Style = $ws.Range("A8").Style

# New value cell (B9) - bold, 14pt, black, Arial font
$ws.Range("B9").Value = "0000-0002-7738-1361"
$font = $ws.Range("B9").Font
$font.Bold = $true
$font.Size = 14
$font.Color = 0
$font.Name = "Arial"

# Row height for the new row (17.4 to match the taller, bigger font)
$ws.Rows.Item(9).RowHeight = 17.4

# Update selection to match the recorded state after the edit
$ws.Range("B9").Select()
